# The "Chart" sheet contains the main GSC export data (Date, Not indexed,
# Indexed, Impressions) starting at row 2. The export was refreshed and no
# longer includes the 2025-10-06 placeholder row (whose Not indexed/Indexed
# columns were blank), so remove that row. Excel shifts every subsequent
# row up by one and drops the now-unused shared strings automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows.Item(2).Delete()
